# "Updated Results with corrected code"
#
# Refresh the per-country results table (Sheet1) with corrected figures
# from the latest model run, and split the old "Other" fuel-demand
# category into "Biogas" + a (smaller) "Other" remainder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Hydrogen: corrected Iron & steel demand (B3); the
# Non-metallic minerals figure (D3) no longer applies and is cleared.
$ws.Range("B3").Value = 902422.6803220445
$ws.Range("D3").ClearContents()

# Row 4 - Methanol: corrected Chemicals demand.
$ws.Range("C4").Value = 10.27572045914962

# Row 5 - Ammonia: corrected Chemicals demand.
$ws.Range("C5").Value = 1831.190154176956

# Row 7 - was "Other", is now specifically "Biogas", with its corrected
# Non-metallic minerals demand.
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 152.7864028263701

# Row 8 (new) - the remaining "Other" demand, broken out separately from
# Biogas above. Match the row-label formatting used by the other
# categories in column A.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 140.182980456325
